$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.868.81"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "'2.317.49"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'270.43"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'94.25"
$ws.Range("E6").Value = "  +7.55%  "

$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("D10").Value = "'44.97"
$ws.Range("E10").Value = "  -1.85%  "

$ws.Range("D11").Value = "'0.0938"
$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").Value = "'8.16"
$ws.Range("E12").Value = "  +7.70%  "

$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").Value = "'2.660.08"
$ws.Range("E14").Value = "  +2.75%  "

$ws.Range("D15").Value = "'15.28"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").Value = "'0.863"
$ws.Range("E16").Value = "  +7.96%  "

$ws.Range("D17").Value = "'2.307.73"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").Value = "'43.825.94"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").Value = "'0.0000107"
$ws.Range("E19").Value = "  +2.28%  "

$ws.Range("D20").Value = "'6.30"
$ws.Range("E20").Value = "  +4.10%  "

$ws.Range("D21").Value = "'71.63"
$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("D22").Value = "'2.30"
$ws.Range("E22").Value = "  -4.49%  "

$ws.Range("D23").Value = "'238.96"
$ws.Range("E23").Value = "  +2.20%  "

$ws.Range("D24").Value = "'9.65"
$ws.Range("E24").Value = "  +9.62%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").Value = "'11.40"
$ws.Range("E26").Value = "  +3.59%  "

$ws.Range("D27").Value = "'2.52"
$ws.Range("E27").Value = "  -2.16%  "

$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("E29").Value = "  -5.01%  "

$ws.Range("D30").Value = "'39.09"
$ws.Range("E30").Value = "  -3.86%  "

$ws.Range("D31").Value = "'22.58"
$ws.Range("E31").Value = "  +8.60%  "

$ws.Range("D32").Value = "'171.96"
$ws.Range("E32").Value = "  -1.85%  "

$ws.Range("D33").Value = "'0.0902"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("E34").Value = "  +3.00%  "

$ws.Range("D35").Value = "'0.127"
$ws.Range("E35").Value = "  +1.95%  "

$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").Value = "'0.0363"
$ws.Range("E37").Value = "  +2.11%  "

$ws.Range("D38").Value = "'4.51"
$ws.Range("E38").Value = "  +1.91%  "

$ws.Range("D39").Value = "'3.49"
$ws.Range("E39").Value = "  +3.19%  "

$ws.Range("E40").Value = "  +14.65%  "

$ws.Range("D41").Value = "'2.31"
$ws.Range("E41").Value = "  +7.39%  "

$ws.Range("D42").Value = "'12.27"
$ws.Range("E42").Value = "  -4.46%  "

$ws.Range("E43").Value = "  +17.12%  "

$ws.Range("D44").Value = "'5.48"
$ws.Range("E44").Value = "  +1.02%  "

$ws.Range("D45").Value = "'61.92"
$ws.Range("E45").Value = "  -5.22%  "

$ws.Range("D46").Value = "'9.00"
$ws.Range("E46").Value = "  +7.02%  "

$ws.Range("E47").Value = "  +3.26%  "

$ws.Range("D48").Value = "'100.54"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'2.536.88"
$ws.Range("E50").Value = "  +2.65%  "

$ws.Range("E51").Value = "  -2.44%  "
